# ajout de la possibilité de réduire le besoin de transport
# Adds a "need_reduction" column (G) to the "Categorie_year" sheet, and
# duplicates the existing 2020 hypotheses block for years 2035 and 2050.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categorie_year")

# --- New column header "need_reduction" with its own (black) font style ---
$ws.Range("G1").Value = "need_reduction"
$ws.Range("G1").Font.Color = 0

# --- Fill the new column for the existing 2020 rows (2-6) with 0 ---
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0

# --- Duplicate the 2020 category block (rows 2-6) for year 2035 (rows 7-11) ---
$ws.Range("A2:A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial() | Out-Null

$ws.Range("B7").Value = 2035
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

$ws.Range("B8").Value = 2035
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 50
$ws.Range("G8").Value = 0

$ws.Range("B9").Value = 2035
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 3.2
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

$ws.Range("B10").Value = 2035
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.61371527777777779
$ws.Range("G10").Value = 0

$ws.Range("B11").Value = 2035
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 245.48611111111109
$ws.Range("G11").Value = 0

# --- Duplicate the 2020 category block (rows 2-6) again for year 2050 (rows 12-16) ---
$ws.Range("A2:A6").Copy() | Out-Null
$ws.Range("A12").PasteSpecial() | Out-Null

$ws.Range("B12").Value = 2050
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 50
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

$ws.Range("B13").Value = 2050
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 50
$ws.Range("G13").Value = 0

$ws.Range("B14").Value = 2050
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 3.2
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

$ws.Range("B15").Value = 2050
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0.61371527777777779
$ws.Range("G15").Value = 0

$ws.Range("B16").Value = 2050
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 245.48611111111109
$ws.Range("G16").Value = 0

# --- Make "Categorie_year" the active sheet/tab, with the same selection Excel left behind ---
$ws.Activate() | Out-Null
$ws.Range("G24").Select() | Out-Null
